# Fruta / hortaliza, semanal
#
# The weekly refresh re-shuffled the existing price-report rows (rows 2-30,
# columns A:R) into a new order. Every row's data still appears somewhere in
# the sheet afterwards -- only the order of rows 2..30 changed (one row, 11,
# and another, 26, happen to land back on themselves).
#
# Strategy: snapshot each data row (A:R) as a Variant array, then write the
# rows back out in their new order, keyed by a row->row map describing which
# original row's data now lives in each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 30

# destinationRow -> sourceRow (which original row's data now occupies it)
$rowMap = @{
    2  = 19
    3  = 27
    4  = 15
    5  = 17
    6  = 20
    7  = 21
    8  = 24
    9  = 22
    10 = 23
    11 = 11
    12 = 2
    13 = 3
    14 = 28
    15 = 18
    16 = 4
    17 = 5
    18 = 6
    19 = 29
    20 = 30
    21 = 13
    22 = 16
    23 = 10
    24 = 25
    25 = 9
    26 = 26
    27 = 14
    28 = 12
    29 = 7
    30 = 8
}

# Snapshot the full A:R block for every original row before any writes happen.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = $ws.Range("A$r" + ":R$r").Value2
}

# Write each destination row from its mapped source row's snapshot.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $rowMap[$r]
    if ($srcRow -ne $r) {
        $ws.Range("A$r" + ":R$r").Value = $snapshot[$srcRow]
    }
}

"Row reshuffle complete"
